$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2237.1936
$ws.Range("I113").Value = 2134.7222
$ws.Range("J113").Value = 2379.077
$ws.Range("K113").Value = 2134.7222
$ws.Range("L113").Value = 2379.077
$ws.Range("M113").Value = 1119.2778
$ws.Range("N113").Value = -8887.077000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 977.2857
$ws.Range("I2").Value = 1068.2
$ws.Range("J2").Value = 750
$ws.Range("K2").Value = 1068.2
$ws.Range("L2").Value = 750
$ws.Range("M2").Value = -955.2
$ws.Range("N2").Value = -976

$ws.Range("H32").Value = 3327.4707
$ws.Range("I32").Value = 2685.6082
$ws.Range("K32").Value = 2685.6082
$ws.Range("M32").Value = -2398.6082

$ws.Range("H43").Value = 10000
$ws.Range("J43").Value = 10000
$ws.Range("L43").Value = 10000
$ws.Range("N43").Value = -10626

$ws.Range("H61").Value = 793.9375
$ws.Range("I61").Value = 804.86664
$ws.Range("J61").Value = 630
$ws.Range("K61").Value = 804.86664
$ws.Range("L61").Value = 630
$ws.Range("M61").Value = -592.86664
$ws.Range("N61").Value = -1054

$ws.Range("H74").Value = 1262.3549
$ws.Range("I74").Value = 826.087
$ws.Range("J74").Value = 2516.625
$ws.Range("K74").Value = 826.087
$ws.Range("L74").Value = 2516.625
$ws.Range("M74").Value = 47.91300000000001
$ws.Range("N74").Value = -4264.625

$ws.Range("H77").Value = 1262.3549
$ws.Range("I77").Value = 826.087
$ws.Range("J77").Value = 2516.625
$ws.Range("K77").Value = 4130.434999999999
$ws.Range("L77").Value = 12583.125
$ws.Range("M77").Value = 237.5650000000005
$ws.Range("N77").Value = -21319.125

$ws.Range("H107").Value = 25000
$ws.Range("J107").Value = 25000
$ws.Range("L107").Value = 25000
$ws.Range("N107").Value = -32680

$ws.Range("H110").Value = 910.5
$ws.Range("I110").Value = 910.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 910.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1134.5
$ws.Range("N110").ClearContents()

$ws.Range("H116").Value = 977.2857
$ws.Range("I116").Value = 1068.2
$ws.Range("J116").Value = 750
$ws.Range("K116").Value = 1068.2
$ws.Range("L116").Value = 750
$ws.Range("M116").Value = 1225.8
$ws.Range("N116").Value = -5338

$ws.Range("H122").Value = 1780.6875
$ws.Range("I122").Value = 1186.5834
$ws.Range("J122").Value = 3563
$ws.Range("K122").Value = 3559.7502
$ws.Range("L122").Value = 10689
$ws.Range("M122").Value = -1109.7502
$ws.Range("N122").Value = -15589

$ws.Range("H136").Value = 793.9375
$ws.Range("I136").Value = 804.86664
$ws.Range("J136").Value = 630
$ws.Range("K136").Value = 2414.59992
$ws.Range("L136").Value = 1890
$ws.Range("M136").Value = 135.4000800000003
$ws.Range("N136").Value = -6990

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 977.2857
$ws.Range("I3").Value = 1068.2
$ws.Range("J3").Value = 750
$ws.Range("K3").Value = 1068.2
$ws.Range("L3").Value = 750
$ws.Range("M3").Value = -954.2
$ws.Range("N3").Value = -978

$ws.Range("H86").Value = 3260.5
$ws.Range("I86").Value = 3171.5
$ws.Range("J86").Value = 3408.8333
$ws.Range("K86").Value = 3171.5
$ws.Range("L86").Value = 3408.8333
$ws.Range("M86").Value = -2048.5
$ws.Range("N86").Value = -5654.8333

$ws.Range("H89").Value = 3260.5
$ws.Range("I89").Value = 3171.5
$ws.Range("J89").Value = 3408.8333
$ws.Range("K89").Value = 15857.5
$ws.Range("L89").Value = 17044.1665
$ws.Range("M89").Value = -10241.5
$ws.Range("N89").Value = -28276.1665

$ws.Range("H107").Value = 1084.68
$ws.Range("I107").Value = 822.44446
$ws.Range("J107").Value = 1759
$ws.Range("K107").Value = 822.44446
$ws.Range("L107").Value = 1759
$ws.Range("M107").Value = 1097.55554
$ws.Range("N107").Value = -5599

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 875069
$ws.Range("J4").Value = 500000
$ws.Range("L4").Value = 1500000
$ws.Range("N4").Value = -1500224

$ws.Range("H113").Value = 2166.3333
$ws.Range("I113").Value = 4194.25
$ws.Range("K113").Value = 12582.75
$ws.Range("M113").Value = -10412.75

$ws.Range("H131").Value = 800.01697
$ws.Range("I131").Value = 366.66666
$ws.Range("J131").Value = 878.02
$ws.Range("K131").Value = 1099.99998
$ws.Range("L131").Value = 2634.06
$ws.Range("M131").Value = 3940.00002
$ws.Range("N131").Value = -12714.06

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2198.75
$ws.Range("I80").Value = 2197.5
$ws.Range("K80").Value = 2197.5
$ws.Range("M80").Value = -1199.5

$ws.Range("H83").Value = 2198.75
$ws.Range("I83").Value = 2197.5
$ws.Range("K83").Value = 10987.5
$ws.Range("M83").Value = -5995.5

$ws.Range("H122").Value = 1651.1852
$ws.Range("I122").Value = 1643.28
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 4929.84
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -2479.84
$ws.Range("N122").Value = -10150

$ws.Range("H126").Value = 56923.11
$ws.Range("I126").Value = 144045.72
$ws.Range("J126").Value = 1481.4546
$ws.Range("K126").Value = 432137.16
$ws.Range("L126").Value = 4444.3638
$ws.Range("M126").Value = -429667.16
$ws.Range("N126").Value = -9384.363799999999

$ws.Range("H132").Value = 3117.0667
$ws.Range("I132").Value = 2827.5386
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 8482.6158
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -5952.6158
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4075
$ws.Range("I40").Value = 3050
$ws.Range("J40").Value = 6125
$ws.Range("K40").Value = 3050
$ws.Range("L40").Value = 6125
$ws.Range("M40").Value = -2914
$ws.Range("N40").Value = -6397

$ws.Range("H139").Value = 57800
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 57800
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 57800
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -68080

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 880
$ws.Range("I107").Value = 880
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2640
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -720
$ws.Range("N107").ClearContents()

$ws.Range("H132").Value = 1649.4193
$ws.Range("I132").Value = 1516.2858
$ws.Range("J132").Value = 1759.0588
$ws.Range("K132").Value = 4548.857400000001
$ws.Range("L132").Value = 5277.1764
$ws.Range("M132").Value = -2018.857400000001
$ws.Range("N132").Value = -10337.1764

$ws.Range("H136").Value = 539.76666
$ws.Range("I136").Value = 269.22726
$ws.Range("J136").Value = 1283.75
$ws.Range("K136").Value = 807.68178
$ws.Range("L136").Value = 3851.25
$ws.Range("M136").Value = 1742.31822
$ws.Range("N136").Value = -8951.25
